$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "prodUIYs"
$ws.Range("B3").Value = "proddjNn"
$ws.Range("B5").Value = "prodLCHk"
